$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for strains that are no longer used (b5145, b5158, b5161, b5164, b5189)
# Delete from bottom to top so earlier row indices stay valid.
$ws.Rows.Item(48).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(4).Delete()

$ws.Range("F73").Select()
